# Weekly fruit/vegetable price update.
# Insert a new week's worth of data (2 rows: "Primera" and "Segunda" quality
# grades) at the top of the data table (row 94), shifting all the existing
# historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 94, pushing the old
# rows 94:176 down to 96:178.
$ws.Rows("94:95").Insert()

# Row 94 - "Primera" quality grade for the new date.
$ws.Cells.Item(94, 1).Value = 11
$ws.Cells.Item(94, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(94, 3).Value = "Bíobío"
$ws.Cells.Item(94, 4).Value = 44512
$ws.Cells.Item(94, 5).Value = 8
$ws.Cells.Item(94, 6).Value = 100112009
$ws.Cells.Item(94, 7).Value = "Acelga"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 200
$ws.Cells.Item(94, 11).Value = 600
$ws.Cells.Item(94, 12).Value = 700
$ws.Cells.Item(94, 13).Value = 650
$ws.Cells.Item(94, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(94, 15).Value = "Región de Ñuble"
$ws.Cells.Item(94, 16).Value = 650
$ws.Cells.Item(94, 17).Value = 1
$ws.Cells.Item(94, 18).Value = "Hortaliza"

# Row 95 - "Segunda" quality grade for the new date.
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value = 44512
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = 100112009
$ws.Cells.Item(95, 7).Value = "Acelga"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Segunda"
$ws.Cells.Item(95, 10).Value = 100
$ws.Cells.Item(95, 11).Value = 500
$ws.Cells.Item(95, 12).Value = 500
$ws.Cells.Item(95, 13).Value = 500
$ws.Cells.Item(95, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(95, 15).Value = "Región de Ñuble"
$ws.Cells.Item(95, 16).Value = 500
$ws.Cells.Item(95, 17).Value = 1
$ws.Cells.Item(95, 18).Value = "Hortaliza"
